# Insert a new price-report row before row 34 (pushing the existing
# rows 34-129 down to 35-130), then populate the new row with the same
# "Madrigal / Primera" Alcachofa record as the row that used to sit at
# position 34, but dated 2023-08-10 (Excel serial 45148).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(34).Insert()

$ws.Range("A34").Value = 5
$ws.Range("B34").Value = "Macroferia Regional de Talca"
$ws.Range("C34").Value = "Maule"
$ws.Range("D34").Value = 45148
$ws.Range("E34").Value = 7
$ws.Range("F34").Value = 100112013
$ws.Range("G34").Value = "Alcachofa"
$ws.Range("H34").Value = "Madrigal"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 12000
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = 12000
$ws.Range("N34").Value = "$/caja 40 unidades"
$ws.Range("O34").Value = "Provincia del Elquí"
$ws.Range("P34").Value = 300
$ws.Range("Q34").Value = 40
$ws.Range("R34").Value = "Hortaliza"
